$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF values (column I) for rows 23 through 52 from 39.46867469879518 to 38.8625
$ws.Range("I23:I52").Value = 38.8625
